$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.873.53'
$ws.Range('E2').Value = '  +2.38%  '
$ws.Range('D3').Value = '3.980.96'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '614.09'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +14.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '163.37'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +10.12%  '
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.998'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('E9').Value = '  +1.80%  '
$ws.Range('E10').Value = '  +1.27%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '54.17'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.77%  '
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.99'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.15%  '
$ws.Range('D14').Value = '4.619.83'
$ws.Range('E14').Value = '  +0.75%  '
$ws.Range('D15').Value = '3.986.79'
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('E16').Value = '  +8.22%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.13'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.02%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '20.53'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('D20').Value = '72.515.38'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '439.67'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.38%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.89'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +15.21%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '96.58'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('E24').Value = '  -3.49%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '14.37'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.74%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.28'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +10.28%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.41'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.14%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.96'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.98%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.53'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.17%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.43'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.75'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '13.95'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +4.57%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '71.98'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +9.86%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '48.16'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -5.76%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '658.08'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.78%  '
$ws.Range('D37').Value = '0.0₃0904'
$ws.Range('E37').Value = '  +10.96%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.440'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.58%  '
$ws.Range('B39').Value = 'ThetaToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.38'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.147'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.16%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  +5.50%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('E44').Value = '  +1.35%  '
$ws.Range('E45').Value = '  +4.10%  '
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.39'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('D49').Value = '2.917.83'
$ws.Range('E49').Value = '  +12.59%  '
$ws.Range('E50').Value = '  +2.68%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.40'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +4.34%  '
